$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Leading apostrophe forces these numeric-looking strings to be stored as
# text (quote-prefixed), preserving exact formatting (e.g. trailing zeros,
# multi-dot "thousand.thousand.decimal" groupings) exactly like the source data.
$ws.Range("D2").Value = '''25.841.20'
$ws.Range("D3").Value = '''1.735.32'
$ws.Range("D5").Value = '''240.74'
$ws.Range("D7").Value = '''0.5191'
$ws.Range("D10").Value = '''1.740.35'
$ws.Range("D12").Value = '''14.93'
$ws.Range("D13").Value = '''0.6396'
$ws.Range("D14").Value = '''4.600'
$ws.Range("D15").Value = '''77.01'
$ws.Range("D18").Value = '''25.876.00'
$ws.Range("D20").Value = '''0.000006755'
$ws.Range("D21").Value = '''1.962.37'
$ws.Range("D22").Value = '''4.265'
$ws.Range("D23").Value = '''8.595'
$ws.Range("D24").Value = '''5.253'
$ws.Range("D25").Value = '''137.57'
$ws.Range("D27").Value = '''15.17'
$ws.Range("D28").Value = '''1.764'
$ws.Range("D29").Value = '''104.87'
$ws.Range("D30").Value = '''3.932'
$ws.Range("D31").Value = '''0.08239'
$ws.Range("D32").Value = '''3.642'
$ws.Range("D33").Value = '''0.04626'
$ws.Range("D34").Value = '''2.643'
$ws.Range("D35").Value = '''0.9852'
$ws.Range("D36").Value = '''0.6167'
$ws.Range("D37").Value = '''2.683'
$ws.Range("D39").Value = '''1.915'
$ws.Range("D41").Value = '''99.88'
$ws.Range("D43").Value = '''0.7444'
$ws.Range("D44").Value = '''4.995'
$ws.Range("D46").Value = '''6.226'
$ws.Range("D48").Value = '''54.78'
$ws.Range("D50").Value = '''7.570'

# --- Volume(1h) % updates (column E) ---
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +4.37%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("E30").Value = '  +5.32%  '
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("E46").Value = '  -0.82%  '
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("E48").Value = '  +2.02%  '
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("E51").Value = '  -1.08%  '
